$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "resume.txt"
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 253
